$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 101.4370242300828
$ws.Cells.Item(3, 2).Value = 150.3499797448662
$ws.Cells.Item(4, 2).Value = 151.9728515229124
$ws.Cells.Item(5, 2).Value = 163.3596629612691
$ws.Cells.Item(6, 2).Value = 176.0488093123393
$ws.Cells.Item(7, 2).Value = 183.4988015505617
$ws.Cells.Item(8, 2).Value = 198.9898690702179
$ws.Cells.Item(9, 2).Value = 211.3453253405451
$ws.Cells.Item(10, 2).Value = 220.8600366436154
$ws.Cells.Item(11, 2).Value = 231.3530538708553
$ws.Cells.Item(12, 2).Value = 240.4283277305466
$ws.Cells.Item(13, 2).Value = 255.9259801349731
$ws.Cells.Item(14, 2).Value = 266.6673492893227
$ws.Cells.Item(15, 2).Value = 277.0108800807276
$ws.Cells.Item(16, 2).Value = 287.9259340340946
$ws.Cells.Item(17, 2).Value = 298.3931317406788
$ws.Cells.Item(18, 2).Value = 308.937849150324
$ws.Cells.Item(19, 2).Value = 325.4175394927131
$ws.Cells.Item(20, 2).Value = 337.7143707405003
$ws.Cells.Item(21, 2).Value = 341.2812461794886
$ws.Cells.Item(22, 2).Value = 341.7449858454843
$ws.Cells.Item(23, 2).Value = 347.3246968335141
$ws.Cells.Item(24, 2).Value = 367.5102459042889
$ws.Cells.Item(25, 2).Value = 373.1642124421396
$ws.Cells.Item(26, 2).Value = 380.6281461016606
$ws.Cells.Item(27, 2).Value = 388.9709587952946
$ws.Cells.Item(28, 2).Value = 415.2880643798964
$ws.Cells.Item(29, 2).Value = 437.4978522305693
$ws.Cells.Item(30, 2).Value = 448.4061001706738
$ws.Cells.Item(31, 2).Value = 439.4886542169731
$ws.Cells.Item(32, 2).Value = 460.0307898036918
$ws.Cells.Item(33, 2).Value = 467.8163474639193
$ws.Cells.Item(34, 2).Value = 481.64120677123
$ws.Cells.Item(35, 2).Value = 497.598626185464
$ws.Cells.Item(36, 2).Value = 503.316305737075
$ws.Cells.Item(37, 2).Value = 527.0153039849024
$ws.Cells.Item(38, 2).Value = 538.3573138003455
$ws.Cells.Item(39, 2).Value = 545.9132275844604
$ws.Cells.Item(40, 2).Value = 563.8887695284817
$ws.Cells.Item(41, 2).Value = 574.812629860202
$ws.Cells.Item(42, 2).Value = 592.8889755424773
$ws.Cells.Item(43, 2).Value = 595.7392959018432
$ws.Cells.Item(44, 2).Value = 605.8916305172444
$ws.Cells.Item(45, 2).Value = 624.8788325562953
$ws.Cells.Item(46, 2).Value = 632.529950172489
$ws.Cells.Item(47, 2).Value = 644.6888067419999
$ws.Cells.Item(48, 2).Value = 650.7045255107772
$ws.Cells.Item(49, 2).Value = 652.7194404791653
$ws.Cells.Item(50, 2).Value = 660.5322201098958
$ws.Cells.Item(51, 2).Value = 671.4525067528568
$ws.Cells.Item(52, 2).Value = 680.4621471419323
$ws.Cells.Item(53, 2).Value = 694.3953568960471
$ws.Cells.Item(54, 2).Value = 698.3846217514247
$ws.Cells.Item(55, 2).Value = 706.1754432635656
$ws.Cells.Item(56, 2).Value = 710.9830896920272
$ws.Cells.Item(57, 2).Value = 720.3687777948928
$ws.Cells.Item(58, 2).Value = 722.0632509819307
$ws.Cells.Item(59, 2).Value = 725.0573360751989
$ws.Cells.Item(60, 2).Value = 732.9991720540273
$ws.Cells.Item(61, 2).Value = 732.7442870247081
$ws.Cells.Item(62, 2).Value = 735.1545484047505
